$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '75.571.90'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +8.51%  '
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.670.12'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +9.93%  '
# Row 4
$ws.Range("E4").Value = '  -0.35%  '
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '187.69'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +13.19%  '
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '587.82'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +4.34%  '
# Row 7
$ws.Range("E7").Value = '  -0.27%  '
# Row 8
$ws.Range("E8").Value = '  +4.28%  '
# Row 9
$ws.Range("E9").Value = '  +14.01%  '
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '2.667.12'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +9.12%  '
# Row 11
$ws.Range("E11").Value = '  +1.64%  '
# Row 12
$ws.Range("E12").Value = '  +7.01%  '
# Row 13
$ws.Range("E13").Value = '  +0.60%  '
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '75.333.67'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +7.91%  '
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.157.40'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +9.21%  '
# Row 16
$ws.Range("E16").Value = '  +5.23%  '
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.53'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +10.75%  '
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.662.86'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +9.31%  '
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '9.24'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +30.31%  '
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.92'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +10.84%  '
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '371.48'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +9.36%  '
# Row 22
$ws.Range("E22").Value = '  +14.74%  '
# Row 23
$ws.Range("E23").Value = '  +5.23%  '
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.27'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +4.06%  '
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.999'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.07%  '
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '69.85'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +5.77%  '
# Row 27
$ws.Range("E27").Value = '  +9.45%  '
# Row 28
$ws.Range("E28").Value = '  +10.40%  '
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.803.39'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +9.02%  '
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.998'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.12%  '
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0₃0948'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +11.87%  '
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.43'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +16.03%  '
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '519.37'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +15.09%  '
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '7.71'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +5.38%  '
# Row 35
$ws.Range("E35").Value = '  +8.73%  '
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.998'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.46%  '
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '163.01'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.26%  '
# Row 38
$ws.Range("E38").Value = '  +6.77%  '
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '19.17'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +5.59%  '
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '19.38'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.48%  '
# Row 41
$ws.Range("E41").Value = '  -0.02%  '
# Row 42
$ws.Range("B42").Value = 'Aave'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '170.07'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +26.63%  '
# Row 43
$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.98'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +14.07%  '
# Row 44
$ws.Range("E44").Value = '  +11.95%  '
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.330'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +8.32%  '
# Row 46
$ws.Range("E46").Value = '  +10.52%  '
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.37'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +13.21%  '
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '39.11'
$ws.Range("D48").Style = "Normal"
# Row 49
$ws.Range("E49").Value = '  +16.63%  '
# Row 50
$ws.Range("E50").Value = '  +7.84%  '
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.534'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +9.56%  '
